# Weekly fruit/vegetable update: a new daily price record was inserted as
# row 85 on the "Vega Modelo de Temuco - Camote" sheet, pushing every
# following record down by one row (old row 149 becomes row 150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 85; everything below (old rows 85-149)
# shifts down to 86-150, which also grows the sheet's used range to R150.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new price record.
$ws.Cells.Item(85, 1).Value = 10
$ws.Cells.Item(85, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(85, 3).Value = "La Araucanía"
$ws.Cells.Item(85, 4).Value = 44907
$ws.Cells.Item(85, 5).Value = 9
$ws.Cells.Item(85, 6).Value = 100114002
$ws.Cells.Item(85, 7).Value = "Camote"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 40
$ws.Cells.Item(85, 11).Value = 24000
$ws.Cells.Item(85, 12).Value = 24000
$ws.Cells.Item(85, 13).Value = 24000
$ws.Cells.Item(85, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(85, 15).Value = "Perú"
$ws.Cells.Item(85, 16).Value = 1200
$ws.Cells.Item(85, 17).Value = 20
$ws.Cells.Item(85, 18).Value = "Hortaliza"
